$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PPL")

# Row 4 (Inventory) - columns B:F
$ws.Range("B4").Value = 361000000.0
$ws.Range("C4").Value = 351000000.0
$ws.Range("D4").Value = 333000000.0
$ws.Range("E4").Value = 320000000.0
$ws.Range("F4").Value = 332000000.0

# Row 13 (Accounts Payable) - columns B:F
$ws.Range("B13").Value = 965000000.0
$ws.Range("C13").Value = 864000000.0
$ws.Range("D13").Value = 804000000.0
$ws.Range("E13").Value = 833000000.0
$ws.Range("F13").Value = 956000000.0

# Row 22 (Long Term Tax Liability (Deferred)) - columns B:F
$ws.Range("B22").Value = 3690000000.0
$ws.Range("C22").Value = 3691000000.0
$ws.Range("D22").Value = 3402000000.0
$ws.Range("E22").Value = 3340000000.0
$ws.Range("F22").Value = 3212000000.0
